# Apply version bump and minor text/punctuation fixes to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Version: 1.0 -> 1.2.5
$ws.Range("D2").Value = "1.2.5"

# Precondition text: fix "usuario" -> "usuário" and add trailing period.
$preconditionCells = @("B8", "B18", "B27", "B37")
foreach ($cell in $preconditionCells) {
    $ws.Range($cell).Value = "O usuário devidamente autenticado e na tela inicial do sistema."
}

# Step 1 description: add trailing period.
$stepCells = @("B10", "B20", "B29", "B39")
foreach ($cell in $stepCells) {
    $ws.Range($cell).Value = "Chefe Acessa a funcionalidade Minha Conta Bancária (menu)."
}

# Expected result for step 2: add trailing period.
$resultCells = @("D11", "D30", "D40")
foreach ($cell in $resultCells) {
    $ws.Range($cell).Value = "SYSTEM Apresenta os campos (banco/agência/conta corrente) alterados."
}

# TC2 informative message: "conta bancários" -> "conta bancária"
$ws.Range("D21").Value = "SYSTEM Exibe mensagens informativas (MSG403 - Informativos sobre a atualização de conta bancária (dados bancários)) para o usuário sobre a manutenção de informações bancárias."
